$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "302.27"
Set-TextValue $ws.Range("E2") "-0.76%"
Set-TextValue $ws.Range("D3") "37.28"
Set-TextValue $ws.Range("E3") "6.86%"
Set-TextValue $ws.Range("D4") "5.004"
Set-TextValue $ws.Range("E4") "-3.63%"
Set-TextValue $ws.Range("D5") "0.07816"
Set-TextValue $ws.Range("E5") "-0.18%"
Set-TextValue $ws.Range("D6") "2.192"
Set-TextValue $ws.Range("E6") "-7.40%"
Set-TextValue $ws.Range("D7") "8.034"
Set-TextValue $ws.Range("E7") "0.08%"
Set-TextValue $ws.Range("D8") "4.038"
Set-TextValue $ws.Range("E8") "2.06%"
Set-TextValue $ws.Range("D9") "0.9127"
Set-TextValue $ws.Range("E9") "-2.19%"
Set-TextValue $ws.Range("D10") "0.09703"
Set-TextValue $ws.Range("E10") "-4.42%"
Set-TextValue $ws.Range("D11") "0.1878"
Set-TextValue $ws.Range("E11") "1.75%"
Set-TextValue $ws.Range("D12") "0.08673"
Set-TextValue $ws.Range("E12") "0.12%"
Set-TextValue $ws.Range("D13") "0.03530"
Set-TextValue $ws.Range("E13") "6.46%"
Set-TextValue $ws.Range("D14") "0.09967"
Set-TextValue $ws.Range("E14") "0.71%"
Set-TextValue $ws.Range("D15") "0.001482"
Set-TextValue $ws.Range("E15") "-0.38%"
Set-TextValue $ws.Range("D16") "0.005705"
Set-TextValue $ws.Range("E16") "-0.23%"
Set-TextValue $ws.Range("D17") "3.458"
Set-TextValue $ws.Range("E17") "-0.33%"
Set-TextValue $ws.Range("D18") "2.082"
Set-TextValue $ws.Range("E18") "-2.01%"
Set-TextValue $ws.Range("E19") "2.26%"
Set-TextValue $ws.Range("D20") "0.1293"
Set-TextValue $ws.Range("E20") "-0.59%"
Set-TextValue $ws.Range("D21") "4.762"
Set-TextValue $ws.Range("E21") "10.46%"
Set-TextValue $ws.Range("D22") "0.2207"
Set-TextValue $ws.Range("E22") "-0.64%"
Set-TextValue $ws.Range("D23") "0.04641"
Set-TextValue $ws.Range("E23") "1.52%"
Set-TextValue $ws.Range("E24") "1.04%"
Set-TextValue $ws.Range("D25") "0.004793"
Set-TextValue $ws.Range("E25") "7.93%"
Set-TextValue $ws.Range("E26") "-7.53%"
Set-TextValue $ws.Range("E27") "39.71%"
Set-TextValue $ws.Range("D39") "0.01757"
Set-TextValue $ws.Range("E39") "-1.21%"
Set-TextValue $ws.Range("D40") "0.04738"
Set-TextValue $ws.Range("E40") "-1.39%"
Set-TextValue $ws.Range("D41") "0.008056"
Set-TextValue $ws.Range("E41") "4.02%"
Set-TextValue $ws.Range("E42") "-1.39%"
Set-TextValue $ws.Range("D43") "0.007684"
Set-TextValue $ws.Range("E43") "10.07%"
Set-TextValue $ws.Range("D44") "0.002221"
Set-TextValue $ws.Range("E44") "0.57%"
Set-TextValue $ws.Range("D45") "0.01045"
Set-TextValue $ws.Range("E45") "10.34%"
Set-TextValue $ws.Range("D46") "0.00006059"
Set-TextValue $ws.Range("E46") "2.56%"
Set-TextValue $ws.Range("E47") "0.25%"
Set-TextValue $ws.Range("D48") "8.126"
Set-TextValue $ws.Range("E48") "197.39%"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.25%"
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.25%"
